$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-17 Wednesday" "2025-12-18 Thursday"

Replace-Text "139÷7=19, 6" "565÷2=282, 1"
Replace-Text "439÷3=146, 1" "682÷4=170, 2"
Replace-Text "714÷8=89, 2" "311÷3=103, 2"
Replace-Text "132÷5=26, 2" "130÷4=32, 2"
Replace-Text "920÷5=184, 0" "521÷2=260, 1"

Replace-Text "650÷5=130, 0" "124÷5=24, 4"
Replace-Text "890÷2=445, 0" "161÷4=40, 1"
Replace-Text "200÷2=100, 0" "791÷5=158, 1"
Replace-Text "593÷2=296, 1" "385÷6=64, 1"
Replace-Text "776÷5=155, 1" "951÷7=135, 6"

Replace-Text "974÷2=487, 0" "652÷3=217, 1"
Replace-Text "868÷2=434, 0" "362÷4=90, 2"
Replace-Text "785÷5=157, 0" "541÷4=135, 1"
Replace-Text "713÷4=178, 1" "491÷5=98, 1"
Replace-Text "828÷3=276, 0" "572÷6=95, 2"

Replace-Text "755÷2=377, 1" "725÷5=145, 0"
Replace-Text "666÷5=133, 1" "293÷6=48, 5"
Replace-Text "194÷2=97, 0" "776÷5=155, 1"
Replace-Text "622÷5=124, 2" "854÷5=170, 4"
Replace-Text "557÷5=111, 2" "267÷5=53, 2"

Replace-Text "119÷8=14, 7" "100÷6=16, 4"
Replace-Text "297÷6=49, 3" "647÷5=129, 2"
Replace-Text "866÷2=433, 0" "767÷6=127, 5"
Replace-Text "178÷6=29, 4" "691÷5=138, 1"
Replace-Text "167÷2=83, 1" "240÷9=26, 6"

"Replacements done"
